$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: tiny floating-point recalculation adjustments to column K
# (centr_eigen_macro) for the affected rows. These are sub-ULP-scale
# differences produced by re-running the eigenvector-centrality
# computation; values are taken verbatim from the target workbook.
$kUpdates = @{
    "2" = 0.7423890598447433
    "3" = 0.6280353788176077
    "4" = 0.5414506112685872
    "5" = 0.7400106338867257
    "6" = 0.6192501947457054
    "7" = 0.8196150948842338
    "8" = 0.6932295827470416
    "10" = 0.7414111553854236
    "11" = 0.7300163799033653
    "12" = 0.5729568581892687
    "13" = 0.7296017761127144
    "15" = 0.7205526401883152
    "17" = 0.7291671489030337
    "19" = 0.7217444379347059
    "20" = 0.7503994041344022
    "22" = 0.7085243743539135
    "24" = 0.6515685568801917
    "25" = 0.7202576681922632
    "27" = 0.7808668713731202
    "28" = 0.7020627893552882
    "31" = 0.677040595518613
    "32" = 0.6221522568492123
    "33" = 0.791112308703898
    "34" = 0.5811148380473606
    "35" = 0.7056051838996015
    "36" = 0.770607464020739
    "37" = 0.6972540736554399
    "38" = 0.6536855299290218
    "39" = 0.6750867163585549
    "40" = 0.6192501947457051
    "41" = 0.7263578664227305
    "42" = 0.7355499095234124
    "43" = 0.3498426195713069
    "45" = 0.6621504322623794
    "46" = 0.7085243743539138
    "48" = 0.6410418120233884
    "50" = 0.5876751203535706
    "52" = 0.5923475419384225
    "54" = 0.7503994041344021
    "55" = 0.503426796616737
    "57" = 0.6461823403930269
    "58" = 0.651744874255984
    "63" = 0.6797318422492199
    "64" = 0.5581775385070359
    "65" = 0.7503994041344021
    "68" = 0.6444444444444446
    "70" = 0.6533878071413061
    "71" = 0.6716815734420243
    "72" = 0.6736031738775605
    "74" = 0.7200562487191552
    "75" = 0.7204348544233163
    "76" = 0.5376470418813519
    "77" = 0.7636558774525605
    "78" = 0.725549621864008
    "79" = 0.6684707170940085
    "80" = 0.7194489972394086
    "81" = 0.6638037066181904
    "82" = 0.6192501947457052
    "83" = 0.6768721242021917
    "84" = 0.6720556872697133
    "85" = 0.6111111111111113
    "86" = 0.5876751203535705
    "87" = 0.46503133108275946
    "89" = 0.7073807095794546
    "90" = 0.6290211899793112
    "91" = 0.7670059672747954
    "94" = 0.8011677544593898
    "95" = 0.7529364663571052
    "96" = 0.6568028603790268
    "97" = 0.7219433793087969
    "98" = 0.7716570891889998
    "99" = 0.7658848514477261
    "100" = 0.634190019720194
    "101" = 0.6952215599621803
    "104" = 0.43249991370841584
    "106" = 0.6543696995514867
    "107" = 0.7051532762746363
    "108" = 0.703183065392293
    "110" = 0.7781497533091724
    "111" = 0.5932868346495678
    "112" = 0.6848446400614068
    "115" = 0.7605622429272167
    "116" = 0.640603345442122
    "117" = 0.5522465746544815
    "118" = 0.7323990412506669
    "119" = 0.7670059672747956
    "122" = 0.8149605674668184
    "123" = 0.6775443814180122
    "124" = 0.6651706006880543
    "125" = 0.5471466884461431
    "126" = 0.6553268404232157
    "127" = 0.5352996816101626
    "128" = 0.7520380747830199
    "129" = 0.796488073217757
    "130" = 0.3337674828751594
    "133" = 0.6908948344366225
    "136" = 0.7280452596986553
    "137" = 0.5242080512939442
    "139" = 0.6389302142172326
    "140" = 0.603124493497912
    "141" = 0.6548896519663272
    "142" = 0.6839268025096327
    "143" = 0.6780697258712605
    "144" = 0.743576418960531
    "145" = 0.7759907622602042
    "146" = 0.723062279606614
    "148" = 0.46873894098876523
    "149" = 0.6363346460488312
    "150" = 0.8167107515573924
    "151" = 0.653387807141306
    "152" = 0.6599824263157813
    "153" = 0.6766248927741227
    "154" = 0.6921410326398417
    "155" = 0.7539578260896126
    "156" = 0.5729944499546029
    "157" = 0.6749707165392153
    "158" = 0.655205145344012
    "159" = 0.5560727419843988
    "160" = 0.6192501947457052
    "162" = 0.6383479422507055
    "163" = 0.6382203233266485
    "164" = 0.6055530991185532
    "165" = 0.7503994041344022
    "166" = 0.7295283117584016
    "167" = 0.5138798723099931
    "168" = 0.5914481721951657
    "170" = 0.5555541941161155
    "172" = 0.6896139328800207
    "173" = 0.6997447656887485
    "174" = 0.7273681331479638
    "175" = 0.7023888300243784
    "176" = 0.5497246451155332
    "178" = 0.7437212953676028
    "179" = 0.7872163097076944
    "180" = 0.7286742057325313
    "181" = 0.667596046705959
    "182" = 0.5876751203535707
    "183" = 0.7107858251581547
    "184" = 0.7750910953322785
    "185" = 0.5921570582603477
    "186" = 0.7085243743539135
    "187" = 0.7653662811812497
    "188" = 0.651186254450461
    "190" = 0.6210125873430767
    "191" = 0.6599557056359688
    "192" = 0.7917764071627923
    "193" = 0.6190775915470379
    "194" = 0.6386654272696458
    "196" = 0.5801949258430951
    "197" = 0.6813106699907031
    "198" = 0.5950426851247089
    "199" = 0.712421535363471
    "200" = 0.7451447138637767
    "201" = 0.6444444444444446
    "204" = 0.5876751203535707
    "206" = 0.668415633175815
    "208" = 0.6975176159222338
    "209" = 0.6431061396231921
    "210" = 0.583906988326783
    "211" = 0.6337316859671387
    "214" = 0.7297654917292115
    "215" = 0.7503994041344023
    "216" = 0.646661226435627
    "217" = 0.6766073699035827
}

foreach ($row in $kUpdates.Keys) {
    $ws.Cells.Item([int]$row, 11).Value = $kUpdates[$row]
}

Write-Host "Updated" $kUpdates.Count "K-column values"
